$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'243.67"
$ws.Cells.Item(3, 4).Value = "'22.92"
$ws.Cells.Item(4, 4).Value = "'5.420"
$ws.Cells.Item(5, 4).Value = "'0.05926"
$ws.Cells.Item(7, 4).Value = "'0.8101"
$ws.Cells.Item(8, 4).Value = "'0.9139"
$ws.Cells.Item(9, 4).Value = "'0.1421"
$ws.Cells.Item(10, 4).Value = "'0.07430"
$ws.Cells.Item(11, 4).Value = "'0.03327"
$ws.Cells.Item(12, 4).Value = "'0.03087"
$ws.Cells.Item(13, 4).Value = "'0.09326"
$ws.Cells.Item(14, 4).Value = "'3.942"
$ws.Cells.Item(15, 4).Value = "'0.001578"
$ws.Cells.Item(16, 4).Value = "'0.04790"
$ws.Cells.Item(18, 4).Value = "'0.005545"
$ws.Cells.Item(19, 4).Value = "'0.004438"
$ws.Cells.Item(20, 4).Value = "'0.0009814"
$ws.Cells.Item(21, 4).Value = "'0.00007807"
$ws.Cells.Item(22, 4).Value = "'3.652"
$ws.Cells.Item(23, 4).Value = "'6.436"
$ws.Cells.Item(24, 4).Value = "'2.150"
$ws.Cells.Item(40, 4).Value = "'0.03889"
$ws.Cells.Item(41, 4).Value = "'0.006209"
$ws.Cells.Item(42, 4).Value = "'0.1065"
$ws.Cells.Item(43, 4).Value = "'0.003003"
$ws.Cells.Item(44, 4).Value = "'0.006538"
$ws.Cells.Item(45, 4).Value = "'0.00005187"
$ws.Cells.Item(49, 4).Value = "'0.002265"
